# Update "想去人数" (F column) figures across the four sheets to reflect
# the latest scrape (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 2798
$ws.Range("F5").Value  = 955
$ws.Range("F7").Value  = 3013
$ws.Range("F8").Value  = 1916
$ws.Range("F9").Value  = 240
$ws.Range("F10").Value = 72
$ws.Range("F11").Value = 2578
$ws.Range("F12").Value = 579
$ws.Range("F13").Value = 275
$ws.Range("F17").Value = 129
$ws.Range("F18").Value = 9604
$ws.Range("F19").Value = 65
$ws.Range("F21").Value = 8
$ws.Range("F22").Value = 7588
$ws.Range("F23").Value = 12126
$ws.Range("F27").Value = 383
$ws.Range("F28").Value = 580
$ws.Range("F29").Value = 2728
$ws.Range("F31").Value = 216
$ws.Range("F32").Value = 2729
$ws.Range("F33").Value = 1174
$ws.Range("F36").Value = 58
$ws.Range("F37").Value = 4564
$ws.Range("F38").Value = 1147
$ws.Range("F40").Value = 365
$ws.Range("F41").Value = 63

# --- Sheet "演出" -------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3

# --- Sheet "本地生活" ----------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 638
$ws.Range("F4").Value = 205

# --- Sheet "全部类型" ----------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 638
$ws.Range("F5").Value  = 2798
$ws.Range("F7").Value  = 3
$ws.Range("F8").Value  = 955
$ws.Range("F11").Value = 3013
$ws.Range("F13").Value = 1916
$ws.Range("F15").Value = 240
$ws.Range("F16").Value = 2578
$ws.Range("F18").Value = 579
$ws.Range("F19").Value = 275
$ws.Range("F22").Value = 129
$ws.Range("F23").Value = 9604
$ws.Range("F25").Value = 8
$ws.Range("F26").Value = 7588
$ws.Range("F27").Value = 12126
$ws.Range("F31").Value = 580
$ws.Range("F33").Value = 2728
$ws.Range("F38").Value = 58
$ws.Range("F39").Value = 4564
